# Update "About" sheet (RMI data refresh: 2019 dollars -> 2020 dollars,
# and refreshed OCCF conversion factor).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# A24: "2018 dollars" -> "2020 dollars"
$ws.Range("A24").Value = "2020 dollars"

# A21: "million 2019 dollars" -> "million 2020 dollars"
$ws.Range("A21").Value = "million 2020 dollars"

# A18: "billion 2019 dollars" -> "billion 2020 dollars"
$ws.Range("A18").Value = "billion 2020 dollars"

# A26: refreshed conversion factor
$ws.Range("A26").Value = 0.88711067149387013

# B29 / B30: update explanatory notes to reference 2020 dollars
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2020 dollar."'
$ws.Range("B30").Value = "2012 dollars are worth more than 2020 dollars, so we need a"

# C1 held a stray date stamp that is no longer used - remove it entirely
# (value + formatting), shrinking the sheet's used range back to column B.
$ws.Range("C1").Clear()

# Leave the cursor where the editor ended up (on the last note cell).
$ws.Range("B31").Select() | Out-Null
